# Generate Report for Handback
# The localization-status report is regenerated: the "e369dccd..." file has
# now also been handed back (in sync with en-US), so both tracked files move
# into "Handed back" status. Rows are re-sorted by source file name
# (e369dccd... before e52524ed...) on every sheet, and the newly-handed-back
# file's "Latest Handback File"/"Latest Handback DateTime" columns are filled
# in on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

$wsOverview.Range("A3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-09 09:45:37"
$wsZh.Range("E2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
$wsZh.Range("F2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-03-09 09:46:26"
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("A3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-09 09:43:11"
$wsZh.Range("E3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
$wsZh.Range("F3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-03-09 09:44:25"
$wsZh.Range("H3").Value = "Include"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.zh-cn.xlf"
    } elseif ($addr -eq '$E$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
    } elseif ($addr -eq '$F$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.zh-cn.xlf"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
    } elseif ($addr -eq '$C$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.zh-cn.xlf"
    } elseif ($addr -eq '$E$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
    } elseif ($addr -eq '$F$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-09 09:45:47"
$wsDe.Range("E2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
$wsDe.Range("F2").Value = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.de-de.xlf"
$wsDe.Range("G2").Value = "2016-03-09 09:46:42"
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("A3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-09 09:43:21"
$wsDe.Range("E3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
$wsDe.Range("F3").Value = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.de-de.xlf"
$wsDe.Range("G3").Value = "2016-03-09 09:44:41"
$wsDe.Range("H3").Value = "Include"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.de-de.xlf"
    } elseif ($addr -eq '$E$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
    } elseif ($addr -eq '$F$2') {
        $h.TextToDisplay = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.de-de.xlf"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
    } elseif ($addr -eq '$C$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.de-de.xlf"
    } elseif ($addr -eq '$E$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
    } elseif ($addr -eq '$F$3') {
        $h.TextToDisplay = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.de-de.xlf"
    }
}
